$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: remove the stray "_GoBack" bookmark from the "Creation of
# solver class" row's second cell, without otherwise changing its text.
# The bookmark cannot be deleted directly through the Bookmarks
# collection (Word hides/guards "_GoBack"), so the row is rebuilt in
# place: a fresh row with identical text is inserted right after it and
# the old (bookmark-bearing) row is then deleted.
$row2 = $t.Rows.Item(2)
$row3 = $t.Rows.Item(3)
$freshRow2 = $t.Rows.Add($row3)
$freshRow2.Cells.Item(1).Range.Text = "Creation of solver class"
$freshRow2.Cells.Item(2).Range.Text = "To keep it separate from the rest of the code to be called on when needed."
$row2.Delete()

# --- Step 2: add the two new design-decision rows just above the
# trailing blank row. Rows.Add(beforeRow) inserts immediately above the
# given row, so insert the lower row first and the upper row second.
$blankRow = $t.Rows.Item($t.Rows.Count)

$swapRow = $t.Rows.Add($blankRow)
$swapRow.Cells.Item(1).Range.Text = "trySwapPiece()"
$swapRow.Cells.Item(2).Range.Text = "We broke this method up, and created 3 new methods, assigning each method with 1 task, this way we can have better cohesion. "

$undoRow = $t.Rows.Add($blankRow)
$undoRow.Cells.Item(1).Range.Text = "Undo and redo"
$undoCellText = "Created in the Game class." + [char]13 + "Undo and Redo are both created as stacks of type move. This way it stores all the move that are made in FILO so you are able to grab the last move made as it will be at the top of the stack. The undo stack stores all the moves made by the user, and redo stores all the moves undone by the function undomove()."
$undoRow.Cells.Item(2).Range.Text = $undoCellText
